$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.110.74'
$ws.Range("E2").Value = '  -1.77%  '
$ws.Range("D3").Value = '1.789.70'
$ws.Range("E3").Value = '  -0.96%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '222.76'
$ws.Range("E5").Value = '  -0.99%  '
$ws.Range("E6").Value = '  -0.65%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '32.35'
$ws.Range("E8").Value = '  -0.50%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.284'
$ws.Range("E9").Value = '  -1.51%  '
$ws.Range("E10").Value = '  +0.02%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0930'
$ws.Range("E11").Value = '  -0.20%  '
$ws.Range("D12").Value = '2.046.65'
$ws.Range("E12").Value = '  -0.93%  '
$ws.Range("D13").Value = '1.797.15'
$ws.Range("E13").Value = '  -0.47%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.95'
$ws.Range("E14").Value = '  -1.63%  '
$ws.Range("E15").Value = '  -2.54%  '
$ws.Range("D16").Value = '34.091.80'
$ws.Range("E17").Value = '  -3.94%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '68.08'
$ws.Range("E18").Value = '  -2.43%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '244.29'
$ws.Range("E19").Value = '  -4.25%  '
$ws.Range("E20").Value = '  -3.32%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.999'
$ws.Range("E21").Value = '  +0.06%  '
$ws.Range("E22").Value = '  -1.17%  '
$ws.Range("E23").Value = '  -4.32%  '
$ws.Range("E24").Value = '  -1.07%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '158.80'
$ws.Range("E25").Value = '  -1.50%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '16.37'
$ws.Range("E26").Value = '  -0.89%  '
$ws.Range("E27").Value = '  -1.66%  '
$ws.Range("E28").Value = '  -2.44%  '
$ws.Range("E29").Value = '  +0.09%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0519'
$ws.Range("E30").Value = '  -3.03%  '
$ws.Range("E31").Value = '  -0.39%  '
$ws.Range("E32").Value = '  -3.81%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.81'
$ws.Range("E34").Value = '  -5.05%  '
$ws.Range("D35").Value = '1.394.63'
$ws.Range("E35").Value = '  -3.76%  '
$ws.Range("E36").Value = '  +1.33%  '
$ws.Range("E37").Value = '  -1.51%  '
$ws.Range("E38").Value = '  -3.90%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '79.74'
$ws.Range("E39").Value = '  -6.82%  '
$ws.Range("E40").Value = '  +0.88%  '
$ws.Range("E41").Value = '  -4.43%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.70'
$ws.Range("E42").Value = '  -3.68%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.17'
$ws.Range("E43").Value = '  +1.41%  '
$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.96'
$ws.Range("E44").Value = '  -2.15%  '
$ws.Range("B45").Value = 'Kaspa'
$ws.Range("C45").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0498'
$ws.Range("E45").Value = '  +0.87%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '107.54'
$ws.Range("E46").Value = '  +1.44%  '
$ws.Range("E47").Value = '  -1.05%  '
$ws.Range("D48").Value = '1.947.08'
$ws.Range("E48").Value = '  -0.50%  '
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '11.99'
$ws.Range("E49").Value = '  -1.19%  '
$ws.Range("B50").Value = 'PaxDollar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.999'
$ws.Range("E50").Value = '  -0.08%  '
$ws.Range("D51").Value = '0.0₆0129'
$ws.Range("E51").Value = '  +1.50%  '
